$wb = $excel.ActiveWorkbook

# Repayment schedule sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new column before column N (shifts old N/O/P -> O/P/Q)
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Update selection on Repayment schedule sheet, activate it
$ws.Range("K19").Select()
$ws.Activate()
